$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 51.66
$ws.Range("I15").Value = 51.66
$ws.Range("K15").Value = 154.98
$ws.Range("M15").Value = 14.02000000000001
$ws.Range("H32").Value = 2299.1
$ws.Range("I32").Value = 923
$ws.Range("J32").Value = 3216.5
$ws.Range("K32").Value = 923
$ws.Range("L32").Value = 3216.5
$ws.Range("M32").Value = -597
$ws.Range("N32").Value = -3868.5
$ws.Range("H58").Value = 8393
$ws.Range("I58").Value = 445.4
$ws.Range("J58").Value = 21639
$ws.Range("K58").Value = 1336.2
$ws.Range("L58").Value = 64917
$ws.Range("M58").Value = -1186.2
$ws.Range("N58").Value = -65217
$ws.Range("H82").Value = 6461.5
$ws.Range("I82").Value = 1323.5
$ws.Range("J82").Value = 11599.5
$ws.Range("K82").Value = 3970.5
$ws.Range("L82").Value = 34798.5
$ws.Range("M82").Value = -3564.5
$ws.Range("N82").Value = -35610.5
$ws.Range("H85").Value = 6461.5
$ws.Range("I85").Value = 1323.5
$ws.Range("J85").Value = 11599.5
$ws.Range("K85").Value = 3970.5
$ws.Range("L85").Value = 34798.5
$ws.Range("M85").Value = -2566.5
$ws.Range("N85").Value = -37606.5
$ws.Range("H87").Value = 20875.334
$ws.Range("J87").Value = 20875.334
$ws.Range("L87").Value = 20875.334
$ws.Range("N87").Value = -23371.334
$ws.Range("H90").Value = 20875.334
$ws.Range("J90").Value = 20875.334
$ws.Range("L90").Value = 62626.00199999999
$ws.Range("N90").Value = -75106.00199999999
$ws.Range("H116").Value = 440715.25
$ws.Range("I116").Value = 668510.0600000001
$ws.Range("J116").Value = 13600
$ws.Range("K116").Value = 668510.0600000001
$ws.Range("L116").Value = 13600
$ws.Range("M116").Value = -665068.0600000001
$ws.Range("N116").Value = -20484
$ws.Range("H129").Value = 1353.8406
$ws.Range("I129").Value = 783
$ws.Range("J129").Value = 1379.7878
$ws.Range("K129").Value = 2349
$ws.Range("L129").Value = 4139.3634
$ws.Range("M129").Value = 2651
$ws.Range("N129").Value = -14139.3634
$ws.Range("H138").Value = 3023.5967
$ws.Range("I138").Value = 1520.0435
$ws.Range("J138").Value = 3910.3076
$ws.Range("K138").Value = 4560.1305
$ws.Range("L138").Value = 11730.9228
$ws.Range("M138").Value = 579.8694999999998
$ws.Range("N138").Value = -22010.9228

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 291.57144
$ws.Range("I2").Value = 274.8
$ws.Range("J2").Value = 333.5
$ws.Range("K2").Value = 274.8
$ws.Range("L2").Value = 333.5
$ws.Range("M2").Value = -161.8
$ws.Range("N2").Value = -559.5
$ws.Range("H32").Value = 3018.204
$ws.Range("I32").Value = 2810.9868
$ws.Range("J32").Value = 3734.0454
$ws.Range("K32").Value = 2810.9868
$ws.Range("L32").Value = 3734.0454
$ws.Range("M32").Value = -2523.9868
$ws.Range("N32").Value = -4308.0454
$ws.Range("H92").Value = 24800
$ws.Range("J92").Value = 24800
$ws.Range("L92").Value = 24800
$ws.Range("N92").Value = -29792
$ws.Range("H116").Value = 291.57144
$ws.Range("I116").Value = 274.8
$ws.Range("J116").Value = 333.5
$ws.Range("K116").Value = 274.8
$ws.Range("L116").Value = 333.5
$ws.Range("M116").Value = 2019.2
$ws.Range("N116").Value = -4921.5
$ws.Range("H131").Value = 41890
$ws.Range("J131").Value = 41890
$ws.Range("L131").Value = 41890
$ws.Range("N131").Value = -51970
$ws.Range("H132").Value = 2591.8635
$ws.Range("I132").Value = 1304.8667
$ws.Range("J132").Value = 5349.7144
$ws.Range("K132").Value = 3914.6001
$ws.Range("L132").Value = 16049.1432
$ws.Range("M132").Value = -1384.6001
$ws.Range("N132").Value = -21109.1432
$ws.Range("H137").Value = 42354.832
$ws.Range("J137").Value = 42354.832
$ws.Range("L137").Value = 42354.832
$ws.Range("N137").Value = -52554.832

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 291.57144
$ws.Range("I3").Value = 274.8
$ws.Range("J3").Value = 333.5
$ws.Range("K3").Value = 274.8
$ws.Range("L3").Value = 333.5
$ws.Range("M3").Value = -160.8
$ws.Range("N3").Value = -561.5
$ws.Range("H134").Value = 3413.78
$ws.Range("I134").Value = 1194.7241
$ws.Range("J134").Value = 6478.1904
$ws.Range("K134").Value = 3584.1723
$ws.Range("L134").Value = 19434.5712
$ws.Range("M134").Value = -1049.1723
$ws.Range("N134").Value = -24504.5712

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1885.8219
$ws.Range("I31").Value = 842.2258
$ws.Range("J31").Value = 2656.0952
$ws.Range("K31").Value = 842.2258
$ws.Range("L31").Value = 2656.0952
$ws.Range("M31").Value = -547.2258
$ws.Range("N31").Value = -3246.0952
$ws.Range("H34").Value = 1885.8219
$ws.Range("I34").Value = 842.2258
$ws.Range("J34").Value = 2656.0952
$ws.Range("K34").Value = 842.2258
$ws.Range("L34").Value = 2656.0952
$ws.Range("M34").Value = -640.2258
$ws.Range("N34").Value = -3060.0952
$ws.Range("H68").Value = 55140.375
$ws.Range("J68").Value = 55140.375
$ws.Range("L68").Value = 55140.375
$ws.Range("N68").Value = -56638.375
$ws.Range("H71").Value = 55140.375
$ws.Range("J71").Value = 55140.375
$ws.Range("L71").Value = 165421.125
$ws.Range("N71").Value = -172909.125
$ws.Range("H86").Value = 2936.3125
$ws.Range("I86").Value = 2552.4614
$ws.Range("K86").Value = 2552.4614
$ws.Range("M86").Value = -1429.4614
$ws.Range("H89").Value = 2936.3125
$ws.Range("I89").Value = 2552.4614
$ws.Range("K89").Value = 12762.307
$ws.Range("M89").Value = -7146.307000000001
$ws.Range("H139").Value = 44518.332
$ws.Range("J139").Value = 44518.332
$ws.Range("L139").Value = 44518.332
$ws.Range("N139").Value = -54798.332

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 290.06668
$ws.Range("I40").Value = 96.5
$ws.Range("K40").Value = 386
$ws.Range("M40").Value = -317
$ws.Range("H68").Value = 1628.6274
$ws.Range("I68").Value = 1480.2
$ws.Range("J68").Value = 1664.8292
$ws.Range("K68").Value = 4440.6
$ws.Range("L68").Value = 4994.487599999999
$ws.Range("M68").Value = -3629.6
$ws.Range("N68").Value = -6616.487599999999
$ws.Range("H71").Value = 1628.6274
$ws.Range("I71").Value = 1480.2
$ws.Range("J71").Value = 1664.8292
$ws.Range("K71").Value = 13321.8
$ws.Range("L71").Value = 14983.4628
$ws.Range("M71").Value = -9265.800000000001
$ws.Range("N71").Value = -23095.4628
$ws.Range("H113").Value = 4167271.8
$ws.Range("I113").Value = 628.94116
$ws.Range("J113").Value = 9615959
$ws.Range("K113").Value = 1886.82348
$ws.Range("L113").Value = 28847877
$ws.Range("M113").Value = 283.17652
$ws.Range("N113").Value = -28852217

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 27890
$ws.Range("J42").Value = 27890
$ws.Range("L42").Value = 27890
$ws.Range("N42").Value = -28860
$ws.Range("H102").Value = 3153.889
$ws.Range("I102").Value = 2221.0952
$ws.Range("J102").Value = 6418.6665
$ws.Range("K102").Value = 2221.0952
$ws.Range("L102").Value = 6418.6665
$ws.Range("M102").Value = -599.0952000000002
$ws.Range("N102").Value = -9662.666499999999
$ws.Range("H115").Value = 27890
$ws.Range("J115").Value = 27890
$ws.Range("L115").Value = 27890
$ws.Range("N115").Value = -30240
$ws.Range("H132").Value = 3383.9697
$ws.Range("I132").Value = 1831.8889
$ws.Range("J132").Value = 5246.467
$ws.Range("K132").Value = 5495.6667
$ws.Range("L132").Value = 15739.401
$ws.Range("M132").Value = -2965.6667
$ws.Range("N132").Value = -20799.401

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1385.0769
$ws.Range("I16").Value = 1091.909
$ws.Range("J16").Value = 2997.5
$ws.Range("K16").Value = 1091.909
$ws.Range("L16").Value = 2997.5
$ws.Range("M16").Value = -921.9090000000001
$ws.Range("N16").Value = -3337.5
$ws.Range("H22").Value = 17930086
$ws.Range("J22").Value = 3100
$ws.Range("L22").Value = 3100
$ws.Range("N22").Value = -3690
$ws.Range("H27").Value = 17930086
$ws.Range("J27").Value = 3100
$ws.Range("L27").Value = 3100
$ws.Range("N27").Value = -3314
$ws.Range("H40").Value = 6150.857
$ws.Range("I40").Value = 5555.636
$ws.Range("J40").Value = 8333.333000000001
$ws.Range("K40").Value = 5555.636
$ws.Range("L40").Value = 8333.333000000001
$ws.Range("M40").Value = -5419.636
$ws.Range("N40").Value = -8605.333000000001
$ws.Range("H122").Value = 6317.6924
$ws.Range("I122").Value = 5613
$ws.Range("J122").Value = 8666.666999999999
$ws.Range("K122").Value = 16839
$ws.Range("L122").Value = 26000.001
$ws.Range("M122").Value = -14389
$ws.Range("N122").Value = -30900.001
$ws.Range("H132").Value = 4051.4412
$ws.Range("I132").Value = 3118.4167
$ws.Range("J132").Value = 6290.7
$ws.Range("K132").Value = 9355.250100000001
$ws.Range("L132").Value = 18872.1
$ws.Range("M132").Value = -6825.250100000001
$ws.Range("N132").Value = -23932.1

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 20419
$ws.Range("J57").Value = 20419
$ws.Range("L57").Value = 20419
$ws.Range("N57").Value = -21927
$ws.Range("H96").Value = 51716676
$ws.Range("I96").Value = 84209176
$ws.Range("J96").Value = 2977915.8
$ws.Range("K96").Value = 84209176
$ws.Range("L96").Value = 2977915.8
$ws.Range("M96").Value = -84207803
$ws.Range("N96").Value = -2980661.8
$ws.Range("H118").Value = 29253.75
$ws.Range("J118").Value = 29253.75
$ws.Range("L118").Value = 29253.75
$ws.Range("N118").Value = -32567.75
$ws.Range("H132").Value = 12348322
$ws.Range("I132").Value = 1889.3334
$ws.Range("K132").Value = 5668.0002
$ws.Range("M132").Value = -3138.0002
$ws.Range("H136").Value = 5058.55
$ws.Range("I136").Value = 2286.7778
$ws.Range("K136").Value = 6860.3334
$ws.Range("M136").Value = -4310.3334
